# Applies the changes described by the commit:
#  - "Casos de Uso" sheet, row 12: Estatus (F12) "En proceso" -> "Hecho"
#    and Q12 (day-9 consumed hours) 0 -> 1 with downstream formulas
#    recalculating automatically.
#  - row 13: Estatus (F13) "Por iniciar" -> "En proceso" and Q13 0 -> 1.
#  - update the active selection / scroll position on the sheet.
#  - re-merge a handful of header merged cells so they end up re-ordered
#    at the end of the mergeCells list (matches how Excel rewrites the
#    list after touching those merges).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")
$ws.Activate()

# --- Estatus column updates -------------------------------------------------
$ws.Range("F12").Value = "Hecho"
$ws.Range("F13").Value = "En proceso"

# --- Consumed hours (Q column) updates --------------------------------------
$ws.Range("Q12").Value = 1
$ws.Range("Q13").Value = 1

# --- Re-touch a subset of the merged header cells ---------------------------
# (reproduces the re-ordering of <mergeCells> seen in the saved file)
$ws.Range("AZ4:BA4").UnMerge()
$ws.Range("AZ4:BA4").Merge()
$ws.Range("AO4:AP4").UnMerge()
$ws.Range("AO4:AP4").Merge()
$ws.Range("AR4:AS4").UnMerge()
$ws.Range("AR4:AS4").Merge()
$ws.Range("AU4:AV4").UnMerge()
$ws.Range("AU4:AV4").Merge()
$ws.Range("AX4:AY4").UnMerge()
$ws.Range("AX4:AY4").Merge()

# --- View / selection state ---------------------------------------------
$ws.Range("R13").Select()
